$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 74 so the current (unedited) row 73 data can be
# relocated there, then overwrite row 73 with the new week's figures.
$ws.Rows.Item(74).Insert()

# Copy row 73 (the original, still-unedited data) down into row 74.
for ($col = 1; $col -le 18; $col++) {
    $src = $ws.Cells.Item(73, $col)
    $dst = $ws.Cells.Item(74, $col)
    $dst.Value2 = $src.Value2
}
# Column D (Fecha) keeps the date number format.
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(73, 4).NumberFormat

# Now update row 73 in place with the new week's data.
$ws.Cells.Item(73, 4).Value2 = 44615   # D73 Fecha
$ws.Cells.Item(73, 11).Value2 = 8000   # K73 Precio minimo
$ws.Cells.Item(73, 12).Value2 = 9000   # L73 Precio maximo
$ws.Cells.Item(73, 13).Value2 = 8500   # M73 Precio promedio ponderado
$ws.Cells.Item(73, 14).Value2 = "$/caja 60 unidades"   # N73 Unidad de comercializacion
$ws.Cells.Item(73, 15).Value2 = "Región del Maule"     # O73 Origen
$ws.Cells.Item(73, 16).Value2 = 142    # P73 Precio $/Kg
$ws.Cells.Item(73, 17).Value2 = 60     # Q73 Kg o Unidades
$ws.Cells.Item(73, 18).Value2 = "Hortaliza"            # R73 Clasificacion
